$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# A new trade entry (2026-02-09, qty 100 @ 36.96) needs to be inserted as the
# newest row at the top of the trade log (row 5), pushing the existing
# entries (rows 5-9) down by one row (to rows 6-10).
#
# Rather than using Rows.Insert() (which stamps the whole inserted row with
# formatting copied from the header row above), shift the data down manually
# by copying each existing row into the row below it, bottom-up, then write
# the new entry into the now-vacated row 5.
for ($r = 9; $r -ge 5; $r--) {
    $dest = $r + 1
    $ws.Range("A$r").Copy($ws.Range("A$dest"))
    $ws.Range("B$r").Copy($ws.Range("B$dest"))
    $ws.Range("C$r").Copy($ws.Range("C$dest"))
    $ws.Range("D$r").Copy($ws.Range("D$dest"))
    $ws.Range("E$r").Copy($ws.Range("E$dest"))
    $ws.Range("F$r").Copy($ws.Range("F$dest"))
    $ws.Range("G$r").Copy($ws.Range("G$dest"))
    $ws.Range("J$dest").Formula = "=Index!`$C`$2"
}

# Populate the new trade entry in row 5
$ws.Range("A5").Value = 46062
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 36.96
$ws.Range("F5").Value = 3696
$ws.Range("G5").Value = "~"
$ws.Range("J5").Formula = "=Index!`$C`$2"
